$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (TC_01) : update Module_Descriptions / TestData columns ---
$ws.Range("B2").Value = "Provide Expected  Home Screen Heading Text"
$ws.Range("C2").Value = "new arrival"

# --- Row 3 (TC_02) : new content, C column becomes a number ---
$ws.Range("B3").Value = "Provide Expected Home Screen Slider Count"
$ws.Range("C3").Value = 3

# --- Row 4 (TC_03) : new content, C column becomes a number ---
$ws.Range("B4").Value = "Provide Expected Home Screen Arrival Count"
$ws.Range("C4").Value = 3

# --- New rows 5-7 : TC_04, TC_05, TC_06 ---
$ws.Range("A5").Value = "TC_04"
$ws.Range("B5").Value = "Provide Expected Product Title"

$ws.Range("A6").Value = "TC_05"
$ws.Range("B6").Value = "Provide Expected Product Description"

$ws.Range("A7").Value = "TC_06"
$ws.Range("B7").Value = "Provide Expected Product Review"

# --- Remove the old per-cell / column style (index 4, empty alignment xf)
# from column A (and the column definition itself) so cells fall back
# to the default "Normal" style, then restore A1's header formatting ---
$ws.Columns(1).ClearFormats()
$ws.Range("A1").HorizontalAlignment = -4108

# --- Update selection to match the new data entry focus ---
$ws.Range("A4:A7").Select()
